$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10775351.19246996
$ws.Range("C2").Value = 2333037.159199236
$ws.Range("D2").Value = 26681052.94840628
$ws.Range("E2").Value = 1138278.49584407
$ws.Range("F2").Value = 7896047.847144105
$ws.Range("G2").Value = 1842394.1173711
$ws.Range("H2").Value = 2144817.388797354
$ws.Range("I2").Value = 10775351.19246996
$ws.Range("J2").Value = 45373051
$ws.Range("L2").Value = 29014090.10760552
$ws.Range("M2").Value = 9034326.342988174
$ws.Range("N2").Value = 3987211.506168454
$ws.Range("O2").Value = 42661.91286283184
$ws.Range("P2").Value = 203894.3726773379
$ws.Range("Q2").Value = 246556.2855401697
$ws.Range("S2").Value = 104322.4875383274
$ws.Range("T2").Value = 104322.4875383274
